$d = $word.ActiveDocument

function Split-RunsInRange {
    param($doc, $rangeStart, $rangeEnd, $runsXml)

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $target = $doc.Range($rangeStart, $rangeEnd)
    $target.InsertXML($pkg)
}

# --- Paragraph "{m:userdoc 'zone1'}" ---
# Originally two runs: "{m" and ":userdoc 'zone1'}".
# Target: four runs: "{" / "m" / ":userdoc 'zone1'" / "}" (last one xml:space="preserve").
$zoneRange = $d.Content
$found = $zoneRange.Find.Execute("{m:userdoc 'zone1'}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $runs1 = "<w:r><w:t>{</w:t></w:r>" +
             "<w:r><w:t>m</w:t></w:r>" +
             "<w:r><w:t>:userdoc 'zone1'</w:t></w:r>" +
             "<w:r><w:t xml:space=`"preserve`">}</w:t></w:r>"
    Split-RunsInRange $d $zoneRange.Start $zoneRange.End $runs1
}

# --- Paragraph "{m:self.name}" ---
# Originally two runs: "{m:self." and "name}".
# Target: "{m:self." stays, "name}" splits into "name" / "}" (last one xml:space="preserve").
$nameRange = $d.Content
$found2 = $nameRange.Find.Execute("name}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $runs2 = "<w:r><w:t>name</w:t></w:r>" +
             "<w:r><w:t xml:space=`"preserve`">}</w:t></w:r>"
    Split-RunsInRange $d $nameRange.Start $nameRange.End $runs2
}
